# Auto-generated Excel COM-interop script applying the cryptos.xlsx refresh
# (commit: 'Updated cryptos list ... with GitHub Actions').
#
# All Price/Volume cells in this sheet are plain text (t="inlineStr"),
# e.g. "242.30" or "  -2.06%  ". Assigning a numeric-looking string
# straight to Range.Value lets Excel auto-coerce it to a real number
# (losing the literal text/trailing zeros), so every write below forces
# a text number-format first, then clears the formatting delta again so
# no stray cell style is left behind (matches the source formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "36.500.94"
Set-TextValue "E2" "  -1.50%  "
Set-TextValue "D3" "2.057.71"
Set-TextValue "E3" "  +0.35%  "
Set-TextValue "E4" "  -0.17%  "
Set-TextValue "D5" "242.34"
Set-TextValue "E5" "  -1.85%  "
Set-TextValue "D6" "0.663"
Set-TextValue "E6" "  +0.29%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "D8" "54.70"
Set-TextValue "E8" "  -4.44%  "
Set-TextValue "D9" "58.61"
Set-TextValue "E9" "  -2.25%  "
Set-TextValue "E10" "  -5.70%  "
Set-TextValue "E11" "  -3.07%  "
Set-TextValue "E12" "  -3.00%  "
Set-TextValue "D13" "0.896"
Set-TextValue "E13" "  -0.42%  "
Set-TextValue "D14" "14.72"
Set-TextValue "E14" "  -6.20%  "
Set-TextValue "D15" "2.359.80"
Set-TextValue "E15" "  +0.13%  "
Set-TextValue "E16" "  -6.68%  "
Set-TextValue "D17" "2.038.77"
Set-TextValue "E17" "  -0.71%  "
Set-TextValue "D18" "36.434.16"
Set-TextValue "E18" "  -1.67%  "
Set-TextValue "D19" "16.71"
Set-TextValue "E19" "  -9.52%  "
Set-TextValue "D20" "72.08"
Set-TextValue "E20" "  -3.34%  "
Set-TextValue "D21" "0.0₃0857"
Set-TextValue "E21" "  -4.67%  "
Set-TextValue "D22" "238.22"
Set-TextValue "E22" "  +0.61%  "
Set-TextValue "D23" "5.24"
Set-TextValue "E23" "  -4.09%  "
Set-TextValue "E24" "  +0.12%  "
Set-TextValue "D25" "2.34"
Set-TextValue "E25" "  -4.78%  "
Set-TextValue "E26" "  -2.03%  "
Set-TextValue "E27" "  -1.95%  "
Set-TextValue "D28" "163.00"
Set-TextValue "E28" "  -4.18%  "
Set-TextValue "D29" "20.19"
Set-TextValue "E29" "  +0.80%  "
Set-TextValue "E30" "  -2.00%  "
Set-TextValue "B31" "ImmutableX"
Set-TextValue "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D31" "1.18"
Set-TextValue "E31" "  +4.12%  "
Set-TextValue "B32" "Filecoin"
Set-TextValue "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "5.03"
Set-TextValue "E32" "  -6.82%  "
Set-TextValue "D33" "4.47"
Set-TextValue "E33" "  -7.26%  "
Set-TextValue "E34" "  -4.10%  "
Set-TextValue "E35" "  -0.06%  "
Set-TextValue "B37" "LidoDAOToken"
Set-TextValue "C37" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D37" "2.19"
Set-TextValue "E37" "  -4.24%  "
Set-TextValue "B38" "Kaspa"
Set-TextValue "C38" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.0825"
Set-TextValue "E38" "  -6.11%  "
Set-TextValue "E39" "  -6.66%  "
Set-TextValue "E40" "  -5.41%  "
Set-TextValue "D41" "0.0214"
Set-TextValue "E41" "  -3.60%  "
Set-TextValue "E42" "  -8.65%  "
Set-TextValue "E43" "  -3.75%  "
Set-TextValue "D44" "93.68"
Set-TextValue "E44" "  -5.04%  "
Set-TextValue "D45" "0.0904"
Set-TextValue "E45" "  -9.29%  "
Set-TextValue "D46" "1.393.82"
Set-TextValue "E46" "  +7.48%  "
Set-TextValue "D47" "15.72"
Set-TextValue "E47" "  -8.05%  "
Set-TextValue "D48" "7.41"
Set-TextValue "E48" "  +8.25%  "
Set-TextValue "E49" "  -0.53%  "
Set-TextValue "E50" "  -4.77%  "
Set-TextValue "D51" "2.247.33"
Set-TextValue "E51" "  +0.24%  "
